$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" values may look numeric (e.g. "7.32"), so force text
# formatting before assigning, then clear the temporary format so the
# cell ends up with no explicit style (matching the source data) while
# keeping the value stored as text, exactly like the original inline strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.814.74"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.77%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.488.02"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.68%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.34"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.87"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.13%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.485.97"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.587"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +5.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.32"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.127"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("E12").Value = "  -0.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.083.91"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.58%  "

$ws.Range("E14").Value = "  -1.60%  "

$ws.Range("E15").Value = "  -1.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.61"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.783.93"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.75%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.467.92"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.17%  "

$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.35"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "393.67"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.29"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.552"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.88"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000126"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.58"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.181"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.14%  "

$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.48"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +8.30%  "

$ws.Range("E31").Value = "  +4.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.07"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.65"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.76"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("E36").Value = "  +2.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.55"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.56"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.99"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.92%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.055.49"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0776"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.40"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.36%  "

$ws.Range("E43").Value = "  -0.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.54"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.89"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.777"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.10"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +11.84%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.13"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.15%  "

$ws.Range("E49").Value = "  +3.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.75"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "312.96"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.84%  "
